$wb = $excel.ActiveWorkbook

# Update values in the "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1829
$ws1.Range("G6").Value = 60
$ws1.Range("G7").Value = 55
$ws1.Range("F12").Value = 5251
$ws1.Range("F13").Value = 68
$ws1.Range("F14").Value = 866
$ws1.Range("F15").Value = 126
$ws1.Range("F16").Value = 2331
$ws1.Range("F18").Value = 39
$ws1.Range("F19").Value = 2173

# Update values in the "全部类型" sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1829
$ws4.Range("G6").Value = 60
$ws4.Range("G7").Value = 55
$ws4.Range("F12").Value = 5251
$ws4.Range("F14").Value = 68
$ws4.Range("F16").Value = 866
$ws4.Range("F17").Value = 126
$ws4.Range("F18").Value = 2331
$ws4.Range("F21").Value = 39
$ws4.Range("F22").Value = 2173
